# Orders.xlsx - "Compeleted core project 2 automation logic"
#
# The automation run's per-row Status/Notes output (columns C & D) is cleared
# out (reset) on every row, the demo product/quantity values are refreshed,
# and two new rows are appended to the bottom of the order list to exercise
# the (now completed) automation logic against an unknown beverage name.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

# --- Row 2 : Chai / -50 -------------------------------------------------
$ws.Range("A2").Value = "Chai"
$ws.Range("B2").Value = -50
$ws.Range("C2").Clear() | Out-Null
$ws.Range("D2").ClearContents() | Out-Null

# --- Row 3 : Ipoh Coffee / 3 --------------------------------------------
$ws.Range("A3").Value = "Ipoh Coffee"
$ws.Range("B3").Value = 3
$ws.Range("C3").Clear() | Out-Null

# --- Row 4 : Outback Lager / 1 ------------------------------------------
$ws.Range("A4").Value = "Outback Lager"
$ws.Range("B4").Value = 1
$ws.Range("C4").Clear() | Out-Null

# --- Row 5 : Guarana Fantastica / 2 --------------------------------------
$ws.Range("A5").Value = "Guarana Fantastica"
$ws.Range("B5").Value = 2
$ws.Range("C5").Clear() | Out-Null

# --- Row 6 : Steeleye Stout / 10 ------------------------------------------
$ws.Range("A6").Value = "Steeleye Stout"
$ws.Range("B6").Value = 10
$ws.Range("C6").Clear() | Out-Null

# --- Row 7 : Laughing Lumberjack Lager / 2 --------------------------------
$ws.Range("A7").Value = "Laughing Lumberjack Lager"
$ws.Range("B7").Value = 2
$ws.Range("C7").Clear() | Out-Null

# --- Row 8 : Ipoh Coff / 4 -------------------------------------------------
$ws.Range("A8").Value = "Ipoh Coff"
$ws.Range("B8").Value = 4
$ws.Range("C8").Clear() | Out-Null

# --- Row 9 : Laughing Lumberjack Lager / 35 --------------------------------
$ws.Range("A9").Value = "Laughing Lumberjack Lager"
$ws.Range("B9").Value = 35
$ws.Range("C9").Clear() | Out-Null
$ws.Range("D9").Clear() | Out-Null

# --- Row 10 : Outback Lager / 2 --------------------------------------------
$ws.Range("A10").Value = "Outback Lager"
$ws.Range("B10").Value = 2
$ws.Range("C10").Clear() | Out-Null
$ws.Range("D10").Clear() | Out-Null

# --- Row 11 : Rhonbrau Klosterbier / 3 --------------------------------------
$ws.Range("A11").Value = "Rhonbrau Klosterbier"
$ws.Range("B11").Value = 3
$ws.Range("C11").Clear() | Out-Null

# --- Row 12 (new) : Not a real beverage lol / 4 -----------------------------
$ws.Rows.Item(12).Insert() | Out-Null
$ws.Range("A12").Value = "Not a real beverage lol"
$ws.Range("B12").Value = 4

# Rows 2/9/10 no longer carry wrapped multi-line notes, so their custom row
# heights collapse back down to the sheet's default.
$ws.Range("2:10").AutoFit() | Out-Null

# Reset the view: scroll back to the top and select the newly added cell.
$ws.Range("B12").Select() | Out-Null

Write-Output "Orders worksheet updated"
